$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137; this shifts the existing rows
# 137..192 down to 138..193 and extends the used range to R193.
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row with the new daily price record.
$ws.Cells.Item(137,1).Value = 1
$ws.Cells.Item(137,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(137,3).Value = "Arica y Parinacota"
$ws.Cells.Item(137,4).Value = 45229
$ws.Cells.Item(137,5).Value = 15
$ws.Cells.Item(137,6).Value = 100112008
$ws.Cells.Item(137,7).Value = "Coliflor"
$ws.Cells.Item(137,8).Value = "Sin especificar"
$ws.Cells.Item(137,9).Value = "Tercera"
$ws.Cells.Item(137,10).Value = 1200
$ws.Cells.Item(137,11).Value = 400
$ws.Cells.Item(137,12).Value = 500
$ws.Cells.Item(137,13).Value = 450
$ws.Cells.Item(137,14).Value = "`$/unidad"
$ws.Cells.Item(137,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(137,16).Value = 450
$ws.Cells.Item(137,17).Value = 1
$ws.Cells.Item(137,18).Value = "Hortaliza"
